{"js": "// The old (Slovenian) text that identifies the target paragraphs. All runs\n// inside each of these paragraphs get collapsed into a single plain run\n// (no run-level formatting) containing the new, translated text.\nconst oldMarker = \"2018: Datumi kampanje za opazovanje\";\nconst newText = \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Gemini: 14. in 23. februar, 14. in 24. marca\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (text.indexOf(oldMarker) !== -1) {\n    // Remove every run in the paragraph (keeping the paragraph mark /\n    // paragraph properties intact), then insert a brand-new run that has\n    // no explicit run properties, just like the target OOXML.\n    paragraph.clear();\n    paragraph.insertText(newText, Word.InsertLocation.start);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The old (Slovenian) text that identifies the target paragraphs. All runs\n# inside each of these paragraphs get replaced by a single plain run (no\n# run-level formatting) containing the new, translated text.\n$oldMarker = \"2018: Datumi kampanje za opazovanje\"\n$newText = \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Gemini: 14. in 23. februar, 14. in 24. marca\"\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $range = $p.Range\n    $text = $range.Text\n    if ($text -like \"*$oldMarker*\") {\n        # Select the paragraph's text but exclude the trailing paragraph\n        # mark so the paragraph itself (and its formatting) is preserved.\n        $range.MoveEnd(1, -1) | Out-Null\n        $range.Delete()\n        $range.InsertAfter($newText)\n    }\n}\n"}
